$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters (G stays 1, H is intentionally empty in every data row)
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
    2 = @(1.02,1.01828226581948,1.023988810243556,1.019570303576058,1.032675055307309,1.029077743253026,1.023491584934031,1.026818630414383,1.022413159273685,1.035479602737922,1.024945060380764)
    3 = @(1.02,1.019206354220505,1.024670923708495,1.020353090682892,1.034851136299598,1.029315555125902,1.024052060315076,1.02730807989642,1.023002063624256,1.037460854153484,1.025506331701136)
    4 = @(1.02,1.019804394387212,1.025111917567211,1.020860085200296,1.036252952206666,1.029467165763846,1.024414229130816,1.027623732581918,1.023382949944778,1.03873627738264,1.025869014838136)
    5 = @(1.02,1.020055832901629,1.025297220643573,1.02107333992435,1.036840811411833,1.029530361229935,1.024566366168211,1.027756181045443,1.023543032732151,1.039270921035129,1.026021367927574)
    6 = @(1.02,1.020098051852938,1.025328328524612,1.021109153016591,1.036939430647859,1.029540940312631,1.024591903685208,1.027778404962443,1.023569908880115,1.039360600384393,1.026046941710775)
    7 = @(1.02,1.019807754032285,1.025114393952212,1.020862934271778,1.036260812919085,1.029468012310651,1.024416262459558,1.027625503353938,1.023385089143066,1.038743427343921,1.025871051054438)
    8 = @(1.02,1.018594546607203,1.02421941221994,1.019834750732507,1.033411787014654,1.029158583841673,1.023681103008589,1.02698426059864,1.022612218251368,1.036150556769322,1.025134847592724)
    9 = @(1.02,1.016457426711641,1.022639432540368,1.018026645188983,1.028342043584558,1.028595871632465,1.022381853045494,1.025846219183335,1.021248991119593,1.031529806652337,1.023833752545656)
    10 = @(1.02,1.015033140805882,1.021584147919493,1.016823738929021,1.024926889155087,1.028208878953126,1.021513118004311,1.025082055894061,1.020339280882006,1.028412522839552,1.02296378380107)
    11 = @(1.02,1.014416513960576,1.021126728254122,1.016303462610997,1.023439234340456,1.028038469811257,1.021136333558272,1.024749859662138,1.019945153998308,1.027053535997131,1.022586464277899)
    12 = @(1.02,1.014187485539689,1.020956750478327,1.016110297601937,1.022885278995918,1.027974743417159,1.020996286021141,1.024626269800599,1.01979872505214,1.026547329608178,1.022446217857202)
    13 = @(1.02,1.014236612277655,1.020993214542761,1.016151728132706,1.02300416715218,1.027988432389959,1.0210263309093,1.024652789191089,1.019830136039683,1.026655977399457,1.022476305412549)
    14 = @(1.02,1.014397582112449,1.021112679304573,1.016287493719312,1.023393472492685,1.028033210929947,1.021124759085141,1.024739647717127,1.019933050803538,1.027011721975261,1.022574873367689)
    15 = @(1.02,1.014496762811515,1.021186275965582,1.016371155145912,1.023633153026431,1.028060743567716,1.02118539159448,1.024793137940852,1.019996455651136,1.027230718861024,1.022635591982145)
    16 = @(1.02,1.015074066326243,1.021614495321329,1.016858280398233,1.025025429497276,1.028220128437186,1.021538110925557,1.025104075027851,1.02036543321322,1.028502517578245,1.022988812215131)
    17 = @(1.02,1.015436219533643,1.021882978669474,1.017163999271384,1.025896363872806,1.028319344691182,1.021759197128122,1.025298766773161,1.020596825055842,1.029297797888668,1.023210212385462)
    18 = @(1.02,1.015647466900574,1.022039534720265,1.01734237684078,1.026403512129042,1.028376942117547,1.021888093521878,1.025412200981854,1.020731771142252,1.029760788115887,1.023339291826882)
    19 = @(1.02,1.015719498424851,1.022092908541788,1.017403208632394,1.026576293132432,1.028396534991445,1.021932033789231,1.0254508577346,1.020777780736433,1.029918507149918,1.023383294494455)
    20 = @(1.02,1.015397362940311,1.021854177682011,1.017131192640895,1.025803009387289,1.02830872805688,1.021735482822868,1.02527789124995,1.020572001059999,1.029212563481475,1.023186464403175)
    21 = @(1.02,1.014350180119355,1.02107750189426,1.016247511686196,1.023278870041672,1.02802003662364,1.021095777011747,1.024714075501693,1.019902745874095,1.026907003526377,1.022545850136427)
    22 = @(1.02,1.013691856848303,1.020588759641656,1.015692420100441,1.021683872207984,1.0278360424472,1.020693030185361,1.024358439697181,1.019481769435922,1.02544918175121,1.022142531363351)
    23 = @(1.02,1.014040838869525,1.020847890693856,1.015986635921456,1.022530180350638,1.027933817348941,1.020906585074091,1.024547077537476,1.019704955002437,1.026222793345425,1.022356389524519)
    24 = @(1.02,1.015414920533319,1.02186719175024,1.017146016369619,1.02584519489508,1.028313526102817,1.021746198480644,1.025287324390506,1.020583218019975,1.029251079969115,1.023197195278414)
    25 = @(1.02,1.017009841080973,1.023048240598505,1.018493644968583,1.029658762880169,1.028743426664544,1.02271819216453,1.026141391793402,1.021601575312083,1.032730718232885,1.02417056930482)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

Write-Output "Updated $($data.Keys.Count) rows"